$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 16 (2015 - Høst): update existing note to include new item
$ws.Range("D16").Value = "2d, 3a-3c"

# Row 14 (2014 - Høst): add new "Ikke lenger pensum" note
$ws.Range("D14").Value = "1f"

# Update selection to D14
$ws.Range("D14").Select()
